$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value to apply. Values are taken verbatim from
# the target OOXML diff (cryptos price/volume refresh commit).
$updates = [ordered]@{
    'D2' = '67.247.24'
    'E2' = '  +1.17%  '
    'D3' = '2.529.52'
    'E3' = '  -1.99%  '
    'D4' = '1.00'
    'E4' = '  -0.06%  '
    'D5' = '590.89'
    'E5' = '  +1.34%  '
    'D6' = '173.79'
    'E6' = '  +4.47%  '
    'E7' = '  -0.04%  '
    'E8' = '  +0.18%  '
    'D9' = '2.528.85'
    'E9' = '  -2.02%  '
    'E10' = '  +0.80%  '
    'E11' = '  +2.03%  '
    'E12' = '  -0.21%  '
    'E13' = '  -3.50%  '
    'E14' = '  -0.49%  '
    'D15' = '2.989.97'
    'E15' = '  -1.91%  '
    'E16' = '  -0.62%  '
    'D17' = '66.996.95'
    'E17' = '  +0.94%  '
    'D18' = '2.530.64'
    'E18' = '  -1.90%  '
    'E19' = '  +4.80%  '
    'D20' = '11.40'
    'E20' = '  -0.20%  '
    'D21' = '354.77'
    'E21' = '  +0.89%  '
    'E22' = '  -0.84%  '
    'E23' = '  +0.50%  '
    'D24' = '2.00'
    'E24' = '  +6.77%  '
    'E25' = '  +0.02%  '
    'D26' = '69.81'
    'E26' = '  +1.59%  '
    'D27' = '9.95'
    'E27' = '  +0.16%  '
    'B28' = 'Binance-PegBSC-USD'
    'C28' = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
    'D28' = '0.998'
    'E28' = '  -0.76%  '
    'B29' = 'WrappedeETH'
    'C29' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D29' = '2.653.81'
    'E29' = '  -2.31%  '
    'D30' = '0.0₃0981'
    'E30' = '  -0.36%  '
    'D31' = '533.51'
    'E31' = '  +0.26%  '
    'D32' = '8.17'
    'E32' = '  +2.12%  '
    'E33' = '  +0.42%  '
    'E34' = '  +0.13%  '
    'E35' = '  -1.03%  '
    'D36' = '1.00'
    'E36' = '  -0.02%  '
    'D37' = '1.47'
    'E37' = '  +0.33%  '
    'D38' = '157.49'
    'E38' = '  +0.33%  '
    'D39' = '18.65'
    'E39' = '  -0.48%  '
    'E40' = '  +0.91%  '
    'D41' = '0.355'
    'E41' = '  -1.39%  '
    'D42' = '1.80'
    'E42' = '  +1.27%  '
    'D43' = '5.14'
    'E43' = '  +0.74%  '
    'E44' = '  -0.02%  '
    'E45' = '  +3.70%  '
    'D46' = '149.06'
    'E46' = '  +0.03%  '
    'E47' = '  -1.39%  '
    'E48' = '  -2.76%  '
    'E49' = '  -0.37%  '
    'D50' = '1.70'
    'E50' = '  -0.36%  '
    'D51' = '0.0758'
    'E51' = '  -0.20%  '
}

foreach ($ref in $updates.Keys) {
    $range = $ws.Range($ref)
    if ($ref.StartsWith("D")) {
        # Price column: force text format so values such as "1.00" or
        # "67.247.24" are preserved verbatim instead of being parsed as numbers.
        $range.NumberFormat = "@"
    }
    $range.Value = $updates[$ref]
}
